$wb = $excel.ActiveWorkbook

# Update the status text everywhere it appears ("Ready for handoff" -> "Handback transform failed")
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("L3").Value = "Handback file name: x2kyg3ye.yxh is different with handoff file name: 0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("L3").Value = "Handback file name: x2kyg3ye.yxh is different with handoff file name: 0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.de-de."
